$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new cell D5 with the dynamic attribute placeholder text
$ws.Range("D5").Value = "{dynamic.name}"

# Update the active selection on the sheet to G24
$ws.Range("G24").Select()
